# Actualización desde MV -datos-
# Append 14 new daily rows (16-09-2021 .. 29-09-2021) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("16-09-2021", 890, 5067, 75),
    @("17-09-2021", 890, 5067, 75),
    @("18-09-2021", 890, 5067, 75),
    @("19-09-2021", 890, 5067, 75),
    @("20-09-2021", 892, 5075, 76),
    @("21-09-2021", 883, 5027, 75),
    @("22-09-2021", 887, 5046, 75),
    @("23-09-2021", 887, 5047, 75),
    @("24-09-2021", 884, 5032, 75),
    @("25-09-2021", 884, 5032, 75),
    @("26-09-2021", 884, 5032, 75),
    @("27-09-2021", 882, 5021, 75),
    @("28-09-2021", 875, 4980, 74),
    @("29-09-2021", 872, 4961, 74)
)

$startRow = 260
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
